$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Plan1")
$ws2 = $wb.Worksheets.Item("Plan2")

# ---------------------------------------------------------------------------
# Plan2 ("Exemplo") gets the same kind of "buteco calc" table as Plan1, but
# for a smaller/different round: 3 drink rows (cerveja, tropeiro simples,
# caipirinha) and 3 people (A, B, C) instead of 4 rows / 8 people.
# ---------------------------------------------------------------------------

# Column widths (mirrors Plan1's D/E/O column widths, applied to Plan2's D/E/J)
$ws2.Columns.Item(4).ColumnWidth = 17
$ws2.Columns.Item(4).Font.Name = $ws2.Columns.Item(4).Font.Name
$ws2.Range("D1").EntireColumn.AutoFit() | Out-Null
$ws2.Columns.Item(4).ColumnWidth = 17
$ws2.Columns.Item(5).ColumnWidth = 12
$ws2.Columns.Item(10).ColumnWidth = 15.7109375

# Header row 4: merged "Quem" banner across F4:H4
$ws2.Range("F4").Value = "Quem"
$ws2.Range("F4:H4").Merge() | Out-Null
$ws2.Range("F4").RowHeight = 15.75

# Row 5 headers
$ws2.Range("C5").Value = "Preco unitario"
$ws2.Range("D5").Value = "Qtd"
$ws2.Range("E5").Value = "Total"
$ws2.Range("F5").Value = "A"
$ws2.Range("G5").Value = "B"
$ws2.Range("H5").Value = "C"
$ws2.Range("I5").Value = "Valor Individual"
$ws2.Range("J5").Value = "Valor por cabeça"

# Row 6: cerveja
$ws2.Range("B6").Value = "cerveja"
$ws2.Range("C6").Value = 10
$ws2.Range("D6").Value = 5
$ws2.Range("E6").Formula = "=D6*C6"
$ws2.Range("F6").Value = "x"
$ws2.Range("G6").Value = "x"
$ws2.Range("H6").Value = "x"
$ws2.Range("I6").Formula = '=COUNTIF(F6:H6,"x")'
$ws2.Range("J6").Formula = "=E6/I6"

# Row 7: tropeiro simples
$ws2.Range("B7").Value = "tropeiro simples"
$ws2.Range("C7").Value = 1
$ws2.Range("D7").Value = 18
$ws2.Range("E7").Formula = "=D7*C7"
$ws2.Range("F7").Value = "x"
$ws2.Range("G7").Value = "x"
$ws2.Range("I7").Formula = '=COUNTIF(F7:H7,"x")'
$ws2.Range("J7").Formula = "=E7/I7"

# Row 8: caipirinha
$ws2.Range("B8").Value = "caipirinha"
$ws2.Range("C8").Value = 3
$ws2.Range("D8").Value = 8
$ws2.Range("E8").Formula = "=D8*C8"
$ws2.Range("G8").Value = "x"
$ws2.Range("H8").Value = "x"
$ws2.Range("I8").Formula = '=COUNTIF(F8:H8,"x")'
$ws2.Range("J8").Formula = "=E8/I8"

# Row 9/10: totals
$ws2.Range("D9").Value = "Total com serviço"
$ws2.Range("E9").Formula = "=SUM(E6:E8)*(1+`$C`$12/100)"
$ws2.Range("D10").Value = "Total sem serviço"
$ws2.Range("E10").Formula = "=SUM(E6:E8)"

# Row 11: blank spacer cell (still styled)
$ws2.Range("E11").Value = $null

# Row 12: service % and per-head totals (with service)
$ws2.Range("B12").Value = "Serviço"
$ws2.Range("C12").Value = 10
$ws2.Range("E12").Value = "Total Individual (com 10%)"
$ws2.Range("F12").Formula = "=(E6/I6+E7/I7)*(1+`$C`$12/100)"
$ws2.Range("G12").Formula = "=(E6/I6+E7/I7+E8/I8)*(1+`$C`$12/100)"
$ws2.Range("H12").Formula = "=(E6/I6+E8/I8)*(1+`$C`$12/100)"
$ws2.Range("I12").Formula = "=SUM(F12:H12)"
$ws2.Range("J12").Value = $null
$ws2.Range("B12").RowHeight = 45

# Row 13: per-head totals (without service)
$ws2.Range("E13").Value = "Total Individual (sem 10%)"
$ws2.Range("F13").Formula = "=(E6/I6+E7/I7)"
$ws2.Range("G13").Formula = "=(E6/I6+E7/I7+E8/I8)"
$ws2.Range("H13").Formula = "=(E6/I6+E8/I8)"
$ws2.Range("I13").Formula = "=SUM(F13:H13)"
$ws2.Range("J13").Value = $null
$ws2.Range("E13").RowHeight = 60

# ---------------------------------------------------------------------------
# Formatting: clone cell-by-cell styles from the analogous Plan1 cells so the
# new Plan2 table looks exactly like Plan1's (same fonts/fills/alignment).
# ---------------------------------------------------------------------------
function CopyFmt($srcRef, $dstRef) {
    $ws1.Range($srcRef).Copy() | Out-Null
    $ws2.Range($dstRef).PasteSpecial(-4122) | Out-Null
}

CopyFmt "F4" "F4:H4"
CopyFmt "E5" "E5"
CopyFmt "C5:D5" "C5:D5"
CopyFmt "F5:H5" "F5:H5"
CopyFmt "N5" "I5"
CopyFmt "O5" "J5"

CopyFmt "B6" "B6:B8"
CopyFmt "C6:D9" "C6:D9"
CopyFmt "E6" "E6:E9"
CopyFmt "F6:M9" "F6:H8"
CopyFmt "N6" "I6:I8"
CopyFmt "O6" "J6:J8"

CopyFmt "D10" "D9:D10"
CopyFmt "E10" "E9:E10"

CopyFmt "E12" "E11"

CopyFmt "B13" "B12"
CopyFmt "C13" "C12"
CopyFmt "E13" "E12"
CopyFmt "F13:N13" "F12:I12"
CopyFmt "O13" "J12"

CopyFmt "E14" "E13"
CopyFmt "F14:N14" "F13:I13"
CopyFmt "O14" "J13"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# View state: Plan2 becomes the active tab with I13 selected; Plan1 loses its
# "tabSelected" flag and keeps its whole used range (A1:O14) selected.
# ---------------------------------------------------------------------------
$ws1.Range("A1:O14").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("I13").Select() | Out-Null
